$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SampleDataPattern")

# --- Adjust the existing "Food:Away:Dinner" row (row 17): ---
# Dinners moved from 3x/week @ $50 to 2x/week @ $75
$ws.Range("C17").Formula = "=-52*2*75"
$ws.Range("E17").Value2 = 2

# --- Insert a new row for "Food:Away:Lunch" right after it, pushing ---
# --- everything from the old row 18 onward down by one row. ---
$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value2 = "Food:Away:Lunch"
$ws.Range("B18").Value2 = "Qdoba,Chipotle,Subway,Jimmy Johns"
$ws.Range("C18").Formula = "=-52*3*15"
$ws.Range("D18").Value2 = "Weekly"
$ws.Range("E18").Value2 = 3
$ws.Range("F18").Value2 = "High"
$ws.Range("G18").Value2 = "High"
$ws.Range("I18").Value2 = "Want"

# --- Expand Table1 (and its AutoFilter) to include the new row. ---
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:I47"))

# --- Update the selection to reflect where the edit was made. ---
$ws.Range("B18").Select()
